$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "26.726.46"
Set-TextValue $ws.Range("E2") "  -1.52%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.594.35"
Set-TextValue $ws.Range("E3") "  -2.25%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "211.67"
Set-TextValue $ws.Range("E5") "  -2.11%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -1.74%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.04%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.249"
Set-TextValue $ws.Range("E8") "  -1.72%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.0618"
Set-TextValue $ws.Range("E9") "  -0.69%  "

# Row 10
Set-TextValue $ws.Range("D10") "19.62"
Set-TextValue $ws.Range("E10") "  -2.41%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0837"
Set-TextValue $ws.Range("E11") "  -1.68%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.818.74"
Set-TextValue $ws.Range("E12") "  -2.22%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.584.32"
Set-TextValue $ws.Range("E13") "  -3.05%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.05"
Set-TextValue $ws.Range("E14") "  -1.59%  "

# Row 15
Set-TextValue $ws.Range("E15") "  -1.89%  "

# Row 16
Set-TextValue $ws.Range("D16") "26.757.71"
Set-TextValue $ws.Range("E16") "  -1.36%  "

# Row 17
Set-TextValue $ws.Range("D17") "63.62"
Set-TextValue $ws.Range("E17") "  -2.58%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.0₃0729"
Set-TextValue $ws.Range("E18") "  -0.17%  "

# Row 19
Set-TextValue $ws.Range("D19") "208.86"
Set-TextValue $ws.Range("E19") "  -2.29%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.73"
Set-TextValue $ws.Range("E21") "  -1.50%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.26"
Set-TextValue $ws.Range("E22") "  -2.65%  "

# Row 23
Set-TextValue $ws.Range("D23") "2.35"
Set-TextValue $ws.Range("E23") "  -5.44%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -2.54%  "

# Row 25
Set-TextValue $ws.Range("D25") "146.56"
Set-TextValue $ws.Range("E25") "  -0.33%  "

# Row 26
Set-TextValue $ws.Range("D26") "7.47"
Set-TextValue $ws.Range("E26") "  +2.29%  "

# Row 27
Set-TextValue $ws.Range("E27") "  +0.05%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.112"
Set-TextValue $ws.Range("E28") "  -4.80%  "

# Row 29
Set-TextValue $ws.Range("D29") "15.33"
Set-TextValue $ws.Range("E29") "  -1.44%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0500"
Set-TextValue $ws.Range("E30") "  -0.74%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -2.12%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -3.23%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.668"
Set-TextValue $ws.Range("E33") "  +23.92%  "

# Row 34
Set-TextValue $ws.Range("E34") "  -2.11%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.311.87"
Set-TextValue $ws.Range("E35") "  -0.44%  "

# Row 36
Set-TextValue $ws.Range("E36") "  -3.47%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -0.92%  "

# Row 38
Set-TextValue $ws.Range("E38") "  -1.19%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.819"
Set-TextValue $ws.Range("E39") "  -2.92%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.01%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.788"
Set-TextValue $ws.Range("E41") "  -2.00%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.18"
Set-TextValue $ws.Range("E42") "  -4.23%  "

# Row 43
Set-TextValue $ws.Range("D43") "5.30"
Set-TextValue $ws.Range("E43") "  +0.80%  "

# Row 44
Set-TextValue $ws.Range("D44") "63.09"
Set-TextValue $ws.Range("E44") "  +0.86%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.730.92"
Set-TextValue $ws.Range("E45") "  -2.09%  "

# Row 46
Set-TextValue $ws.Range("D46") "89.03"
Set-TextValue $ws.Range("E46") "  -1.99%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.62"
Set-TextValue $ws.Range("E47") "  +1.34%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.821"
Set-TextValue $ws.Range("E48") "  +1.82%  "

# Row 49
Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0509"
Set-TextValue $ws.Range("E49") "  -0.76%  "

# Row 50
Set-TextValue $ws.Range("B50") "Algorand"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.0980"
Set-TextValue $ws.Range("E50") "  +3.41%  "

# Row 51
Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.47"
Set-TextValue $ws.Range("E51") "  -1.36%  "
